$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.355.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "'1.786.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'226.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'32.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.67%  "
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").Value = "'0.0688"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").Value = "'0.0945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "'2.044.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'11.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.77%  "
$ws.Range("D14").Value = "'1.781.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "'34.390.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "'68.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").Value = "'245.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "'11.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "'168.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").Value = "'7.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.25%  "
$ws.Range("D27").Value = "'16.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'4.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.21%  "
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'3.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("D35").Value = "'1.412.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").Value = "'2.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.71%  "
$ws.Range("E37").Value = "  +4.87%  "
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "'84.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.39%  "
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").Value = "'14.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("D47").Value = "'6.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "'1.946.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").Value = "'105.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -0.94%  "
